# ExcelCard.xlsx - fillCells: push text into the template's placeholder
# cells and apply the matching "Item" look (bold Arial, centered header,
# left-aligned body) cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$boldCenterItem = "Жирный Центр Item"
$boldLeftItem   = "Жирный лево Item"
$finalLeft      = "final left"

# --- Row 4: title, centered, Arial 11 bold ---
$c = $ws.Range("B4")
$c.Value = $boldCenterItem
$c.Font.Name = "Arial"
$c.Font.Size = 11
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

# --- Row 5: left label, Arial 10 bold ---
$c = $ws.Range("B5")
$c.Value = $boldLeftItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

# --- Row 6: left label + centered value ---
$c = $ws.Range("B6")
$c.Value = $boldLeftItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

$c = $ws.Range("C6")
$c.Value = $boldCenterItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

# --- Row 7: left label + final-left value ---
$c = $ws.Range("B7")
$c.Value = $boldLeftItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

$c = $ws.Range("D7")
$c.Value = $finalLeft
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

# --- Row 8: left label + centered value + final-left value ---
$c = $ws.Range("B8")
$c.Value = $boldLeftItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

$c = $ws.Range("C8")
$c.Value = $boldCenterItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

$c = $ws.Range("D8")
$c.Value = $finalLeft
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

# --- Row 9: left label + final-left value ---
$c = $ws.Range("B9")
$c.Value = $boldLeftItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

$c = $ws.Range("D9")
$c.Value = $finalLeft
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

# --- Row 10: final-left value only ---
$c = $ws.Range("C10")
$c.Value = $finalLeft
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

# --- Row 11: left label + final-left value ---
$c = $ws.Range("B11")
$c.Value = $boldLeftItem
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108

$c = $ws.Range("C11")
$c.Value = $finalLeft
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.Font.Bold = $true
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4108
